$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column D (rows 3-10) with the letters a-h, turning them into shared strings
$values = @("a", "b", "c", "d", "e", "f", "g", "h")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Match the paper/orientation page-setup the author left on the sheet
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait

# Leave the selection where it ends up after typing through D10 and pressing Enter
[void]$ws.Range("D11").Select()
